$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("area_pop_sum")

# Fix capitalization of "Population" -> "population"
$ws.Range("A3").Value = "population"

# Add new "density" row
$ws.Range("A4").Value = "density"
$ws.Range("B4").Value = 2185.317890879674
